# Atualização de bases das ligas, do dia: 10-06-2024 às 21:53
#
# Two pairs of consecutive match rows had their data (everything except
# column A, the running sequence number) swapped between each other:
#   - Row 169 <-> Row 170
#   - Row 176 <-> Row 177

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($row1, $row2) {
    $addr1 = "B" + $row1 + ":AD" + $row1
    $addr2 = "B" + $row2 + ":AD" + $row2

    $range1 = $ws.Range($addr1)
    $range2 = $ws.Range($addr2)

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value2 = $values2
    $range2.Value2 = $values1
}

Swap-RowData 169 170
Swap-RowData 176 177
